$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2037735849056604
$ws.Range("C2").Value = 0.5320754716981132
$ws.Range("J2").Value = 0.03773584905660377
$ws.Range("P2").Value = 0.1320754716981132
$ws.Range("S2").Value = 0.09433962264150944
$ws.Range("B3").Value = 0.02684563758389262
$ws.Range("C3").Value = 0.02684563758389262
$ws.Range("J3").Value = 0.08053691275167785
$ws.Range("P3").Value = 0.5906040268456376
$ws.Range("S3").Value = 0.2751677852348993
$ws.Range("J4").Value = 0.1212121212121212
$ws.Range("P4").Value = 0.6060606060606061
$ws.Range("S4").Value = 0.2727272727272727
$ws.Range("B6").Value = 0.09248554913294797
$ws.Range("D6").Value = 0.01734104046242774
$ws.Range("F6").Value = 0.04046242774566474
$ws.Range("J6").Value = 0.1849710982658959
$ws.Range("O6").Value = 0.02312138728323699
$ws.Range("Q6").Value = 0.115606936416185
$ws.Range("R6").Value = 0.05780346820809248
$ws.Range("S6").Value = 0.4682080924855491
$ws.Range("B7").Value = 0.1590909090909091
$ws.Range("D7").Value = 0.02840909090909091
$ws.Range("F7").Value = 0.05113636363636364
$ws.Range("J7").Value = 0.1306818181818182
$ws.Range("O7").Value = 0.005681818181818182
$ws.Range("Q7").Value = 0.1363636363636364
$ws.Range("R7").Value = 0.07386363636363637
$ws.Range("S7").Value = 0.4147727272727273
$ws.Range("B8").Value = 0.06609195402298851
$ws.Range("D8").Value = 0.005747126436781609
$ws.Range("F8").Value = 0.04310344827586207
$ws.Range("J8").Value = 0.09770114942528736
$ws.Range("O8").Value = 0.02298850574712644
$ws.Range("Q8").Value = 0.235632183908046
$ws.Range("R8").Value = 0.1120689655172414
$ws.Range("S8").Value = 0.4166666666666667
$ws.Range("B9").Value = 0.1257142857142857
$ws.Range("D9").Value = 0.005714285714285714
$ws.Range("E9").Value = 0.005714285714285714
$ws.Range("F9").Value = 0.06857142857142857
$ws.Range("J9").Value = 0.1428571428571428
$ws.Range("O9").Value = 0.02285714285714286
$ws.Range("Q9").Value = 0.1771428571428571
$ws.Range("R9").Value = 0.09714285714285714
$ws.Range("S9").Value = 0.3542857142857143
$ws.Range("B10").Value = 0.1039301310043668
$ws.Range("D10").Value = 0.02183406113537118
$ws.Range("E10").Value = 0.0008733624454148472
$ws.Range("F10").Value = 0.07423580786026202
$ws.Range("J10").Value = 0.1301310043668122
$ws.Range("O10").Value = 0.02008733624454149
$ws.Range("Q10").Value = 0.2043668122270742
$ws.Range("R10").Value = 0.09082969432314411
$ws.Range("S10").Value = 0.3537117903930131
$ws.Range("G11").Value = 0.1052631578947368
$ws.Range("J11").Value = 0.1052631578947368
$ws.Range("K11").Value = 0.1614035087719298
$ws.Range("L11").Value = 0.6175438596491228
$ws.Range("S11").Value = 0.01052631578947368
$ws.Range("G12").Value = 0.7287234042553191
$ws.Range("J12").Value = 0.1861702127659574
$ws.Range("K12").Value = 0.01063829787234043
$ws.Range("L12").Value = 0.02659574468085106
$ws.Range("S12").Value = 0.04787234042553191
$ws.Range("G13").Value = 0.6
$ws.Range("J13").Value = 0.3333333333333333
$ws.Range("S13").Value = 0.06666666666666667
$ws.Range("F15").Value = 0.02183406113537118
$ws.Range("H15").Value = 0.1266375545851528
$ws.Range("I15").Value = 0.09170305676855896
$ws.Range("J15").Value = 0.3755458515283843
$ws.Range("K15").Value = 0.07860262008733625
$ws.Range("M15").Value = 0.008733624454148471
$ws.Range("O15").Value = 0.08296943231441048
$ws.Range("S15").Value = 0.2139737991266376
$ws.Range("F16").Value = 0.0145985401459854
$ws.Range("H16").Value = 0.1605839416058394
$ws.Range("I16").Value = 0.1094890510948905
$ws.Range("J16").Value = 0.4087591240875912
$ws.Range("K16").Value = 0.1313868613138686
$ws.Range("M16").Value = 0.0218978102189781
$ws.Range("O16").Value = 0.0583941605839416
$ws.Range("S16").Value = 0.0948905109489051
$ws.Range("F17").Value = 0.005141388174807198
$ws.Range("H17").Value = 0.1748071979434447
$ws.Range("I17").Value = 0.07197943444730077
$ws.Range("J17").Value = 0.4190231362467866
$ws.Range("K17").Value = 0.1208226221079692
$ws.Range("M17").Value = 0.02313624678663239
$ws.Range("N17").Value = 0.002570694087403599
$ws.Range("O17").Value = 0.07712082262210797
$ws.Range("S17").Value = 0.1053984575835476
$ws.Range("F18").Value = 0.01092896174863388
$ws.Range("H18").Value = 0.1912568306010929
$ws.Range("I18").Value = 0.07103825136612021
$ws.Range("J18").Value = 0.3497267759562842
$ws.Range("K18").Value = 0.1256830601092896
$ws.Range("M18").Value = 0.02185792349726776
$ws.Range("O18").Value = 0.1311475409836066
$ws.Range("S18").Value = 0.09836065573770492
$ws.Range("F19").Value = 0.01293900184842884
$ws.Range("H19").Value = 0.1829944547134935
$ws.Range("I19").Value = 0.08964879852125693
$ws.Range("J19").Value = 0.3909426987060998
$ws.Range("K19").Value = 0.1173752310536044
$ws.Range("M19").Value = 0.0009242144177449168
$ws.Range("N19").Value = 0.0009242144177449168
$ws.Range("O19").Value = 0.07578558225508318
$ws.Range("S19").Value = 0.1146025878003697
